$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Kanban board updates ---
# "Homepage System" task is replaced by two new tasks; "Dashboard System" takes its
# spot in the Not Started column, "Search Book System" is added as a new Not Started task.
$ws.Range("A4").Value = "Dashboard System"
$ws.Range("A7").Value = "Search Book System"

# "View a Book Instance Page" has been initiated: move it from Not Started to Doing.
$ws.Range("B4").Value = "View a Book Instance Page"

# "Add Books System" has been finished: move it from Doing to Done.
$ws.Range("C17").Value = "Add Books System"

# Update the active cell/selection to reflect where the author was last working.
[void]$ws.Range("D13").Select()
